$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.641.11'
$ws.Range('E2').Value = '  -1.86%  '
$ws.Range('D3').Value = '2.890.68'
$ws.Range('E3').Value = '  -2.19%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '570.44'
$ws.Range('E5').Value = '  -4.02%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '143.87'
$ws.Range('E6').Value = '  -1.43%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.502'
$ws.Range('E8').Value = '  -0.94%  '
$ws.Range('D9').Value = '2.890.37'
$ws.Range('E9').Value = '  -2.18%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.65'
$ws.Range('E10').Value = '  -8.13%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.149'
$ws.Range('E11').Value = '  -1.38%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.432'
$ws.Range('E12').Value = '  -2.68%  '
$ws.Range('E13').Value = '  -1.65%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '31.98'
$ws.Range('E14').Value = '  -3.40%  '
$ws.Range('E15').Value = '  -0.78%  '
$ws.Range('D16').Value = '3.370.43'
$ws.Range('E16').Value = '  -2.31%  '
$ws.Range('D17').Value = '61.632.40'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '6.61'
$ws.Range('E18').Value = '  -1.71%  '
$ws.Range('D19').Value = '2.903.24'
$ws.Range('E19').Value = '  -1.91%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '434.70'
$ws.Range('E20').Value = '  -1.58%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.18'
$ws.Range('E21').Value = '  -2.34%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.654'
$ws.Range('E22').Value = '  -2.29%  '
$ws.Range('E23').Value = '  -2.77%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '78.96'
$ws.Range('E24').Value = '  -3.10%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '11.92'
$ws.Range('E25').Value = '  +0.48%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '10.16'
$ws.Range('E26').Value = '  -9.85%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('E28').Value = '  -4.54%  '
$ws.Range('E29').Value = '  +10.03%  '
$ws.Range('E30').Value = '  -2.34%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.51'
$ws.Range('E32').Value = '  -3.98%  '
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.106'
$ws.Range('E34').Value = '  -2.77%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '25.56'
$ws.Range('E35').Value = '  -3.62%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.956'
$ws.Range('E36').Value = '  -3.81%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.43'
$ws.Range('E37').Value = '  -3.97%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '48.99'
$ws.Range('E38').Value = '  -1.05%  '
$ws.Range('E40').Value = '  -4.59%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.116'
$ws.Range('E41').Value = '  -0.97%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '8.25'
$ws.Range('E42').Value = '  -3.16%  '
$ws.Range('E43').Value = '  -4.26%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '38.74'
$ws.Range('E44').Value = '  -5.47%  '
$ws.Range('D45').Value = '2.686.63'
$ws.Range('E45').Value = '  -1.77%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '133.61'
$ws.Range('E46').Value = '  -0.81%  '
$ws.Range('E47').Value = '  -1.87%  '
$ws.Range('E48').Value = '  +0.04%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '337.66'
$ws.Range('E49').Value = '  -7.73%  '
$ws.Range('E50').Value = '  -2.15%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '21.65'
$ws.Range('E51').Value = '  -5.55%  '
